$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.419475
$ws.Range("H2").Value = 4.258425
$ws.Range("I2").Value = 0.1541931834006784
$ws.Range("J2").Value = 0.1541931834006784
$ws.Range("M2").Value = 3.063353333333333
$ws.Range("N2").Value = 9.190059999999999
$ws.Range("O2").Value = 0.1884019917097105
$ws.Range("P2").Value = 0.1884019917097105
$ws.Range("Q2").Value = 4.348353472833333
$ws.Range("R2").Value = 39.13518125549999
$ws.Range("S2").Value = 0.02905030286074848
$ws.Range("T2").Value = 0.02905030286074848
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.419475
$ws.Range("H3").Value = 4.258425
$ws.Range("I3").Value = 0.1541931834006784
$ws.Range("J3").Value = 0.1541931834006784
$ws.Range("O3").Value = 0.4156086771445645
$ws.Range("P3").Value = 0.4156086771445645
$ws.Range("Q3").Value = 9.59232658955
$ws.Range("R3").Value = 86.33093930594998
$ws.Range("S3").Value = 0.06408402497786517
$ws.Range("T3").Value = 0.06408402497786515
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.419475
$ws.Range("H4").Value = 4.258425
$ws.Range("I4").Value = 0.1541931834006784
$ws.Range("J4").Value = 0.1541931834006784
$ws.Range("M4").Value = 3.493414666666666
$ws.Range("N4").Value = 10.480244
$ws.Range("O4").Value = 0.214851572590793
$ws.Range("P4").Value = 0.214851572590793
$ws.Range("Q4").Value = 4.958814783966666
$ws.Range("R4").Value = 44.62933305569999
$ws.Range("S4").Value = 0.03312864793641632
$ws.Range("T4").Value = 0.03312864793641631
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.419475
$ws.Range("H5").Value = 4.258425
$ws.Range("I5").Value = 0.1541931834006784
$ws.Range("J5").Value = 0.1541931834006784
$ws.Range("M5").Value = 2.945239333333333
$ws.Range("N5").Value = 8.835718
$ws.Range("O5").Value = 0.1811377585549322
$ws.Range("P5").Value = 0.1811377585549322
$ws.Range("Q5").Value = 4.180693602683333
$ws.Range("R5").Value = 37.62624242415
$ws.Range("S5").Value = 0.02793020762564846
$ws.Range("T5").Value = 0.02793020762564846
$ws.Range("I6").Value = 0.3984988340349546
$ws.Range("J6").Value = 0.3984988340349546
$ws.Range("M6").Value = 3.063353333333333
$ws.Range("N6").Value = 9.190059999999999
$ws.Range("O6").Value = 0.1884019917097105
$ws.Range("P6").Value = 0.1884019917097105
$ws.Range("Q6").Value = 11.23794029463111
$ws.Range("R6").Value = 101.14146265168
$ws.Range("S6").Value = 0.0750779740261828
$ws.Range("T6").Value = 0.0750779740261828
$ws.Range("I7").Value = 0.3984988340349546
$ws.Range("J7").Value = 0.3984988340349546
$ws.Range("O7").Value = 0.4156086771445645
$ws.Range("P7").Value = 0.4156086771445645
$ws.Range("S7").Value = 0.1656195732569188
$ws.Range("T7").Value = 0.1656195732569188
$ws.Range("I8").Value = 0.3984988340349546
$ws.Range("J8").Value = 0.3984988340349546
$ws.Range("M8").Value = 3.493414666666666
$ws.Range("N8").Value = 10.480244
$ws.Range("O8").Value = 0.214851572590793
$ws.Range("P8").Value = 0.214851572590793
$ws.Range("Q8").Value = 12.81562430987022
$ws.Range("R8").Value = 115.340618788832
$ws.Range("S8").Value = 0.0856181011680074
$ws.Range("T8").Value = 0.0856181011680074
$ws.Range("I9").Value = 0.3984988340349546
$ws.Range("J9").Value = 0.3984988340349546
$ws.Range("M9").Value = 2.945239333333333
$ws.Range("N9").Value = 8.835718
$ws.Range("O9").Value = 0.1811377585549322
$ws.Range("P9").Value = 0.1811377585549322
$ws.Range("Q9").Value = 10.80463798323378
$ws.Range("R9").Value = 97.241741849104
$ws.Range("S9").Value = 0.07218318558384558
$ws.Range("T9").Value = 0.07218318558384558
$ws.Range("G10").Value = 4.049549666666667
$ws.Range("H10").Value = 12.148649
$ws.Range("I10").Value = 0.4398900681184871
$ws.Range("J10").Value = 0.439890068118487
$ws.Range("M10").Value = 3.063353333333333
$ws.Range("N10").Value = 9.190059999999999
$ws.Range("O10").Value = 0.1884019917097105
$ws.Range("P10").Value = 0.1884019917097105
$ws.Range("Q10").Value = 12.40520146988222
$ws.Range("R10").Value = 111.64681322894
$ws.Range("S10").Value = 0.08287616496684318
$ws.Range("T10").Value = 0.08287616496684316
$ws.Range("G11").Value = 4.049549666666667
$ws.Range("H11").Value = 12.148649
$ws.Range("I11").Value = 0.4398900681184871
$ws.Range("J11").Value = 0.439890068118487
$ws.Range("O11").Value = 0.4156086771445645
$ws.Range("P11").Value = 0.4156086771445645
$ws.Range("Q11").Value = 27.36547170134733
$ws.Range("R11").Value = 246.289245312126
$ws.Range("S11").Value = 0.1828221292997567
$ws.Range("T11").Value = 0.1828221292997567
$ws.Range("G12").Value = 4.049549666666667
$ws.Range("H12").Value = 12.148649
$ws.Range("I12").Value = 0.4398900681184871
$ws.Range("J12").Value = 0.439890068118487
$ws.Range("M12").Value = 3.493414666666666
$ws.Range("N12").Value = 10.480244
$ws.Range("O12").Value = 0.214851572590793
$ws.Range("P12").Value = 0.214851572590793
$ws.Range("Q12").Value = 14.14675619892844
$ws.Range("R12").Value = 127.320805790356
$ws.Range("S12").Value = 0.09451107290232801
$ws.Range("T12").Value = 0.09451107290232799
$ws.Range("G13").Value = 4.049549666666667
$ws.Range("H13").Value = 12.148649
$ws.Range("I13").Value = 0.4398900681184871
$ws.Range("J13").Value = 0.439890068118487
$ws.Range("M13").Value = 2.945239333333333
$ws.Range("N13").Value = 8.835718
$ws.Range("O13").Value = 0.1811377585549322
$ws.Range("P13").Value = 0.1811377585549322
$ws.Range("Q13").Value = 11.92689296055356
$ws.Range("R13").Value = 107.342036644982
$ws.Range("S13").Value = 0.07968070094955917
$ws.Range("T13").Value = 0.07968070094955916
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.068288
$ws.Range("H14").Value = 0.204864
$ws.Range("I14").Value = 0.00741791444588001
$ws.Range("J14").Value = 0.007417914445880009
$ws.Range("M14").Value = 3.063353333333333
$ws.Range("N14").Value = 9.190059999999999
$ws.Range("O14").Value = 0.1884019917097105
$ws.Range("P14").Value = 0.1884019917097105
$ws.Range("Q14").Value = 0.2091902724266667
$ws.Range("R14").Value = 1.88271245184
$ws.Range("S14").Value = 0.001397549855936027
$ws.Range("T14").Value = 0.001397549855936027
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.068288
$ws.Range("H15").Value = 0.204864
$ws.Range("I15").Value = 0.00741791444588001
$ws.Range("J15").Value = 0.007417914445880009
$ws.Range("O15").Value = 0.4156086771445645
$ws.Range("P15").Value = 0.4156086771445645
$ws.Range("Q15").Value = 0.461466949504
$ws.Range("R15").Value = 4.153202545536
$ws.Range("S15").Value = 0.003082949610023746
$ws.Range("T15").Value = 0.003082949610023746
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.068288
$ws.Range("H16").Value = 0.204864
$ws.Range("I16").Value = 0.00741791444588001
$ws.Range("J16").Value = 0.007417914445880009
$ws.Range("M16").Value = 3.493414666666666
$ws.Range("N16").Value = 10.480244
$ws.Range("O16").Value = 0.214851572590793
$ws.Range("P16").Value = 0.214851572590793
$ws.Range("Q16").Value = 0.2385583007573333
$ws.Range("R16").Value = 2.147024706816
$ws.Range("S16").Value = 0.001593750584041281
$ws.Range("T16").Value = 0.001593750584041281
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.068288
$ws.Range("H17").Value = 0.204864
$ws.Range("I17").Value = 0.00741791444588001
$ws.Range("J17").Value = 0.007417914445880009
$ws.Range("M17").Value = 2.945239333333333
$ws.Range("N17").Value = 8.835718
$ws.Range("O17").Value = 0.1811377585549322
$ws.Range("P17").Value = 0.1811377585549322
$ws.Range("Q17").Value = 0.2011245035946667
$ws.Range("R17").Value = 1.810120532352
$ws.Range("S17").Value = 0.001343664395878957
$ws.Range("T17").Value = 0.001343664395878957
